$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The commit ("Error Calculations and Plots") removes two rows that were
# re-added by mistake in the "missing data" export (the "RM 232" row, and
# the "SC 92" row a couple rows below it), and fixes up three cells whose
# presence/absence of a value changed between the two exports.

# Delete the "RM 232" row (row 26). Rows below shift up by one.
$ws.Rows.Item(26).Delete()

# After that shift, the "SC 92" row (originally row 28) is now row 27.
# Delete it too; rows below shift up by one more.
$ws.Rows.Item(27).Delete()

# "SC 5" is now row 26 - its column D value, previously missing, is filled in.
$ws.Range("D26").Value = -13.8

# "SC 101" is now row 27 - its column D value becomes missing.
$ws.Range("D27").ClearContents()

# "SC 232" is now row 33 (the new last row) - its column F value,
# previously missing, is filled in.
$ws.Range("F33").Value = 17.53
